# "Added proper player movement"
# Each player token on the Clue board occupies a "current location" cell
# which is drawn with a highlighted fill color and a value ending in "*"
# (e.g. "P*", "L*", "s*", "K*", "b*", "c*", "R*", "A*", "B*"), while the
# rest of that player's cells just show the plain room/letter value with
# the normal board fill (style 1).
#
# This edit advances every player one space: the old "current location"
# cell reverts to its plain value/style, and the new "current location"
# cell receives the highlighted value/style that used to belong to the
# old cell. In other words, for each player the two involved cells swap
# their (value, style) pairs.
#
# B2 additionally had a stray/incorrect fill (style 10, matching the
# Library's highlight color) despite holding the plain "P" value; that is
# corrected to the normal board style (style 1) without altering its value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

$scratch = $ws.Range("AZ100")

function Swap-CellFormats {
    param(
        [string]$Ref1,
        [string]$Ref2
    )

    $c1 = $ws.Range($Ref1)
    $c2 = $ws.Range($Ref2)

    # Stash c1's (value, format) in the scratch cell.
    $c1.Copy() | Out-Null
    $scratch.PasteSpecial($xlPasteFormats) | Out-Null
    $scratch.Value = $c1.Value2

    # Move c2's (value, format) into c1.
    $c2.Copy() | Out-Null
    $c1.PasteSpecial($xlPasteFormats) | Out-Null
    $c1.Value = $c2.Value2

    # Move the stashed original c1 (value, format) into c2.
    $scratch.Copy() | Out-Null
    $c2.PasteSpecial($xlPasteFormats) | Out-Null
    $c2.Value = $scratch.Value2

    $scratch.Clear() | Out-Null
    $ws.Application.CutCopyMode = $false
}

# Each player moves one step: swap the "current location" highlighted
# cell with the plain cell it moves into.
Swap-CellFormats "D4" "D5"     # P  (Miss Scarlet-style pawn): D4 <-> D5
Swap-CellFormats "L4" "L6"     # L  pawn: L4 <-> L6
Swap-CellFormats "T4" "S2"     # C  pawn: T4 <-> S2
Swap-CellFormats "D10" "B13"   # s  pawn: D10 <-> B13
Swap-CellFormats "T11" "T13"   # K  pawn: T11 <-> T13
Swap-CellFormats "L21" "M22"   # b  pawn: L21 <-> M22
Swap-CellFormats "T21" "S22"   # c  pawn: T21 <-> S22
Swap-CellFormats "D26" "D29"   # R  pawn: D26 <-> D29
Swap-CellFormats "T28" "V30"   # A  pawn: T28 <-> V30
Swap-CellFormats "L29" "M31"   # B  pawn: L29 <-> M31

# Fix B2's stray highlight color (value stays "P").
$fmtSrc = $ws.Range("C2")   # normal board style (s=1)
$fmtSrc.Copy() | Out-Null
$ws.Range("B2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Application.CutCopyMode = $false

# Update the saved selection to match the author's final cursor position.
$ws.Range("S9").Select() | Out-Null
